$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Terms")

# Insert a new column before column B (shifts old B..E to C..F)
$ws.Range("B1").EntireColumn.Insert()

# New "Composition" column header (write B2 first so "Composition" gets
# registered in the shared string table before "[Composition]")
$ws.Range("B2").Value = "Composition"
$ws.Range("B1").Value = "[Composition]"

# New data point for the composition row (row 7: s;t)
$ws.Range("B7").Value = "s;t"

# Activate the Terms sheet and select D21 to match the author's final view state
$ws.Activate()
$ws.Range("D21").Select()
